# Correção das notas do fórum para matc65 em 2021.2
# Zera os valores das colunas B:J (visualizações diárias, total_views e nota_view)
# para todas as linhas de dados (linhas 2 a 50), mantendo a coluna A (matricula)
# e o cabeçalho (linha 1) intactos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
